$wb = $excel.ActiveWorkbook
Write-Output "sheets:"
foreach ($s in $wb.Worksheets) { Write-Output $s.Name }
